# Applies the "automatic update of files" edit described by the diff:
# Rows 2-5 (data rows) are refreshed with new values coming from an upstream
# data source. Rows are effectively re-ordered (2->3, 3->5, 4->2, 5->4) and
# the Id (A) / Taxonsorteringsordning (B) values are refreshed to new numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (new content, previously held by row 4, with refreshed A/B ids)
$ws.Range("A2").Value = 111188106
$ws.Range("B2").Value = 78107
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 6453
$ws.Range("F2").Value = "Vedskivlav"
$ws.Range("G2").Value = "Hertelidea botryosa"
$ws.Range("H2").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q2").Value = 367050.1089280251
$ws.Range("R2").Value = 6876649.905198708
$ws.Range("S2").Value = 4
$ws.Range("Z2").Value = "16:29"
$ws.Range("AB2").Value = "16:29"
$ws.Range("AC2").Value = ""

# Row 3 (new content, previously held by row 2, with refreshed A/B ids)
$ws.Range("A3").Value = 111188165
$ws.Range("B3").Value = 90666
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 4364
$ws.Range("F3").Value = "Dropptaggsvamp"
$ws.Range("G3").Value = "Hydnellum ferrugineum"
$ws.Range("H3").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q3").Value = 367058.2797417908
$ws.Range("R3").Value = 6876642.060615195
$ws.Range("S3").Value = 4
$ws.Range("Z3").Value = "16:33"
$ws.Range("AB3").Value = "16:33"
$ws.Range("AC3").Value = ""

# Row 4 (new content, previously held by row 5, with refreshed A/B ids)
$ws.Range("A4").Value = 111191235
$ws.Range("B4").Value = 77515
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("Q4").Value = 367076.7989791847
$ws.Range("R4").Value = 6876754.204911136
$ws.Range("S4").Value = 5
$ws.Range("Z4").Value = "18:08"
$ws.Range("AB4").Value = "18:08"
$ws.Range("AC4").Value = "Växer på gran i fuktig skogsmiljö. Kontinuitetsskog"

# Row 5 (new content, previously held by row 3, with refreshed A/B ids)
$ws.Range("A5").Value = 111191051
$ws.Range("B5").Value = 89423
$ws.Range("E5").Value = 5432
$ws.Range("F5").Value = "Granticka"
$ws.Range("G5").Value = "Porodaedalea chrysoloma"
$ws.Range("H5").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q5").Value = 367080.0096928239
$ws.Range("R5").Value = 6876739.970800492
$ws.Range("S5").Value = 4
$ws.Range("Z5").Value = "17:59"
$ws.Range("AB5").Value = "17:59"
$ws.Range("AC5").Value = "Växer på gammal gran i fuktig skogsmiljö."
